# Applies the "final set of supplementary for Chapter 4" edit:
#  - A1 title text: "respectively" -> "respectively," (comma added)
#  - Column A group labels renamed to short codes (MCL14-BTB, MCL21-BTB, WIA20-BTB, OGR25-BTB)
#  - Column B gene-set labels renamed (13-gene set / 17-gene set / 30-gene set)
#  - Rows 2-14 height changed from 20.1 to 24.95
#  - G3:H5 number formatting reverted to General (no more scientific style)
#  - Several P / P.adj. values in G6:H14 recalculated

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A1 title: insert a comma after "respectively" while preserving the
#     existing rich-text run formatting as closely as possible. ---
$titleRange = $ws.Range("A1")
$fullText = $titleRange.Characters().Text
$idx = $fullText.IndexOf("respectively")
if ($idx -ge 0) {
    $insertPos = $idx + "respectively".Length + 1   # 1-based position right after the word
    $titleRange.Characters($insertPos, 0).Text = ","
}

# Re-apply the original run formatting (the plain .Text assignment above
# collapses rich-text runs, so restore each run's font explicitly).
$titleRange.Characters(1, 27).Font.Bold = $true
$titleRange.Characters(1, 27).Font.Italic = $false
$titleRange.Characters(1, 27).Font.Size = 14
$titleRange.Characters(1, 27).Font.Name = "Aptos Narrow"

$titleRange.Characters(28, 24).Font.Bold = $false
$titleRange.Characters(28, 24).Font.Italic = $false
$titleRange.Characters(28, 24).Font.Size = 14
$titleRange.Characters(28, 24).Font.Name = "Aptos Narrow"

$titleRange.Characters(52, 1).Font.Bold = $false
$titleRange.Characters(52, 1).Font.Italic = $true
$titleRange.Characters(52, 1).Font.Size = 14
$titleRange.Characters(52, 1).Font.Name = "Aptos Narrow"

$titleRange.Characters(53, 28).Font.Bold = $false
$titleRange.Characters(53, 28).Font.Italic = $false
$titleRange.Characters(53, 28).Font.Size = 14
$titleRange.Characters(53, 28).Font.Name = "Aptos Narrow"

$titleRange.Characters(81, 105).Font.Bold = $false
$titleRange.Characters(81, 105).Font.Italic = $false
$titleRange.Characters(81, 105).Font.Size = 14
$titleRange.Characters(81, 105).Font.Name = "Aptos Narrow"

# --- Column A: rename study labels to short codes ---
$ws.Range("A3").Value = "MCL14-BTB"
$ws.Range("A4").Value = "MCL14-BTB"
$ws.Range("A5").Value = "MCL14-BTB"

$ws.Range("A6").Value = "MCL21-BTB"
$ws.Range("A7").Value = "MCL21-BTB"
$ws.Range("A8").Value = "MCL21-BTB"

$ws.Range("A9").Value = "WIA20-BTB"
$ws.Range("A10").Value = "WIA20-BTB"
$ws.Range("A11").Value = "WIA20-BTB"

$ws.Range("A12").Value = "OGR25-BTB"
$ws.Range("A13").Value = "OGR25-BTB"
$ws.Range("A14").Value = "OGR25-BTB"

# --- Column B: rename Pass 1/Pass 2/Combined to gene-set sizes ---
foreach ($r in @(3, 6, 9, 12)) {
    $ws.Cells.Item($r, 2).Value = "13-gene set"
}
foreach ($r in @(4, 7, 10, 13)) {
    $ws.Cells.Item($r, 2).Value = "17-gene set"
}
foreach ($r in @(5, 8, 11, 14)) {
    $ws.Cells.Item($r, 2).Value = "30-gene set"
}

# --- Row heights 2-14: 20.1 -> 24.95 ---
$ws.Range("A2:A14").Rows.RowHeight = 24.95

# --- G3:H5 lose the scientific-notation style (back to the default "Normal" style) ---
$ws.Range("G3:H5").Style = "Normal"

# --- Updated P / P.adj. values ---
$ws.Range("H6").Value  = [double]"2.55428571428571E-6"

$ws.Range("G7").Value  = [double]"2.9699999999999999E-6"
$ws.Range("H7").Value  = [double]"4.4549999999999997E-6"

$ws.Range("G8").Value  = [double]"5.9499999999999998E-6"
$ws.Range("H8").Value  = [double]"7.9333333333333308E-6"

$ws.Range("H9").Value  = [double]"2.88E-8"

$ws.Range("G10").Value = [double]"1.2E-8"
$ws.Range("H10").Value = [double]"2.88E-8"

$ws.Range("G11").Value = [double]"3.2600000000000001E-8"
$ws.Range("H11").Value = [double]"6.5200000000000001E-8"

$ws.Range("G12").Value = [double]"3.8200000000000002E-19"
$ws.Range("H12").Value = [double]"1.5280000000000001E-18"

$ws.Range("G13").Value = [double]"2.3800000000000002E-22"
$ws.Range("H13").Value = [double]"1.4279999999999999E-21"

$ws.Range("G14").Value = [double]"2.3800000000000002E-22"
$ws.Range("H14").Value = [double]"1.4279999999999999E-21"

Write-Output "Edit applied."
